$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                                       $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output "WARNING: not found -> $find"
    }
}

# Objective line: "... Intern in Spring or Summer of 2017" -> "... Intern in the summer of 2017"
Replace-Text "Intern in Spring or Summer of " "Intern in the summer of "

# Projects section: Holodesk -> Healthnet
Replace-Text "Holodesk " "Healthnet "

# Projects section: September 2016 -> August 2016 (date range for the Healthnet project)
Replace-Text "September 2016 " "August 2016 "

# Bullet: Member of a team building ... -> Leader of a team building ...
Replace-Text "Member of a team building" "Leader of a team building"

# Bullet: touch screen display using a projector and Kinect -> health tracking and patient management system in Django
Replace-Text " a touch screen display using a projector and Kinect" " a health tracking and patient management system in Django"

# Bullet: Writing code to improve touch detection using Arduino and ADXL345 accelerometers
#          -> Coordinating all team activity for a four-person development team
Replace-Text "Writing code to improve touch detection using Arduino and ADXL345 accelerometers" "Coordinating all team activity for a four-person development team"

# Extracurricular bullet: many fields of interest -> various technical subjects
Replace-Text "y about many fields of interest" "y about various technical subjects"

# Skills - Experienced with: remove C, add Django + Processing at the end
Replace-Text " Java, Python, Ruby, C, HTML, Arduino, Git, JavaScript" " Java, Python, Ruby, HTML, Arduino, Git, JavaScript, Django, Processing"

# Skills - Exposed to: drop leading Processing, add C (keep "Exposed to:" bold run untouched)
Replace-Text "Processing, SQL, Android, Rails, Clojure" "C, SQL, Android, Rails, Clojure"
